$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (changed) date column C for every data row (2-463)
#    from serial 45184 to serial 45186.
$ws.Range("C2:C463").Value = 45186

# 2) Append a display-text second argument to the HYPERLINK() formulas in
#    columns S, T, V, W, X, Y for every row that has them (rows 2-12), using
#    the "Beteckning" identifier from column A of that same row.
$cols = @("S", "T", "V", "W", "X", "Y")
for ($row = 2; $row -le 12; $row++) {
    $ident = $ws.Range("A$row").Value2
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$row")
        $formula = $cell.Formula
        if ($formula -and $formula.Length -gt 0) {
            # Strip the trailing ')' and append the friendly text argument.
            $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $ident + '")'
            $cell.Formula = $newFormula
        }
    }
}
